$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete rows (old items 8-15, sheet rows 9-16).
# The remaining rows 1-8 will be updated in place below, and the
# used range / dimension will shrink to A1:G8 automatically.
$ws.Rows("9:16").Delete()

# Row 2 (item 1): BAV21A0F1C -> BAV24G0I1C, new description, new sizes/qty
$ws.Range("B2").Value = "BAV24G0I1C"
$ws.Range("C2").Value = "BALL VALVE W/INTEGRAL WELDED 2 NIPPLES, FB, FLOATING BALL, API 608, API 598, A105, CL 800, SW W/2 PE NIPPLES, MNF STD, SS316 BALL, SS316 STEM, 2 OR 3 PCS SPLIT BODY SIDE ENTRY, API 607, LO"
$ws.Range("D2").Value = "0,75"
$ws.Range("E2").Value = "5,00"

# Row 3 (item 2): BAV21A0F1C -> BAV24G0I1C, new description, new sizes/qty, note CSO
$ws.Range("B3").Value = "BAV24G0I1C"
$ws.Range("C3").Value = "BALL VALVE W/INTEGRAL WELDED 2 NIPPLES, FB, FLOATING BALL, API 608, API 598, A105, CL 800, SW W/2 PE NIPPLES, MNF STD, SS316 BALL, SS316 STEM, 2 OR 3 PCS SPLIT BODY SIDE ENTRY, API 607, LO"
$ws.Range("D3").Value = "0,75"
$ws.Range("E3").Value = "2,00"
$ws.Range("G3").Value = "CSO"

# Row 4 (item 3): BAV21A0I1C -> BAV24G0I1C, new description, new sizes/qty
$ws.Range("B4").Value = "BAV24G0I1C"
$ws.Range("C4").Value = "BALL VALVE W/INTEGRAL WELDED 2 NIPPLES, FB, FLOATING BALL, API 608, API 598, A105, CL 800, SW W/2 PE NIPPLES, MNF STD, SS316 BALL, SS316 STEM, 2 OR 3 PCS SPLIT BODY SIDE ENTRY, API 607, LO"
$ws.Range("D4").Value = "1,00"
$ws.Range("E4").Value = "5,00"

# Row 5 (item 4): BAV21A0I1C -> BAV24G0I1C, new description, new size, note CSO
$ws.Range("B5").Value = "BAV24G0I1C"
$ws.Range("C5").Value = "BALL VALVE W/INTEGRAL WELDED 2 NIPPLES, FB, FLOATING BALL, API 608, API 598, A105, CL 800, SW W/2 PE NIPPLES, MNF STD, SS316 BALL, SS316 STEM, 2 OR 3 PCS SPLIT BODY SIDE ENTRY, API 607, LO"
$ws.Range("D5").Value = "1,00"
$ws.Range("G5").Value = "CSO"

# Row 6 (item 5): BAV24G0I1C -> CKV22A0B2B, swing check valve description, new sizes/qty
$ws.Range("B6").Value = "CKV22A0B2B"
$ws.Range("C6").Value = "SWING CHECK VALVE FL, API 594, API 598, A216 GR.WCB, CL 300, INST HORIZ/VERT, RF, B16.5, BOLTED COVER, SPW SS304/GRAPH, RENEWABLE SEATS, TRIM #8"
$ws.Range("D6").Value = "10,00"
$ws.Range("E6").Value = "2,00"

# Row 7 (item 6): BAV24G0I1C -> GAV22A0B2B, gate valve description, new sizes/qty, note cleared
$ws.Range("B7").Value = "GAV22A0B2B"
$ws.Range("C7").Value = "GATE VALVE FL, API 600, API 598, A216 GR.WCB, CL 300, RF, B16.5, BB, SPW SS304/GRAPH, PKG GRAPH, TRIM #8, RENEWABLE SEATS, FLEXIBLE WEDGE, STEM OS&Y/RSNRO, GO"
$ws.Range("D7").Value = "10,00"
$ws.Range("E7").Value = "1,00"
$ws.Range("G7").Value = "-"

# Row 8 (item 7): BAV24G0I1C -> MFV22A0I2I, mono flange description, new sizes/qty
$ws.Range("B8").Value = "MFV22A0I2I"
$ws.Range("C8").Value = "INTEGRAL MONO FLANGE DBB NEEDLE MULTI-VALVE, EEMUA 182, A105, CL 300, RF/NPTF, B16.5 AND B1.20.1, BB, SPW SS304/GRAPH, PKG GRAPH; SS316 STEM, SEATS&STEM TIP, S, SWIVEL NEEDLE, STEM OS&Y/RSRO, T-HANDLE"
$ws.Range("D8").Value = "0,75"
$ws.Range("E8").Value = "4,00"

Write-Output "edit applied"
